$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# J2 holds a zero-padded code ("001") - force text so leading zero survives
$ws.Range("J2").NumberFormat = "@"
$ws.Range("J2").Value = "001"

$ws.Range("N2").Value = "2019-12-31 00:00:00"

$ws.Range("O2").Value = 412955424.54
$ws.Range("P2").Value = 82368939.66
$ws.Range("Q2").Value = 158098549.49
$ws.Range("R2").Value = 12.9478978471
$ws.Range("S2").Value = 34371169.64
$ws.Range("T2").Value = -29.0027734993
$ws.Range("U2").Value = 37455156.32
$ws.Range("V2").Value = 6.5946544741
$ws.Range("W2").Value = 48576852.82
$ws.Range("X2").Value = 25378202.62
$ws.Range("Y2").Value = -19.5874017978
$ws.Range("Z2").Value = 5899594.94
$ws.Range("AA2").Value = 35.0016980261
$ws.Range("AB2").Value = 364378571.72
$ws.Range("AC2").Value = 12.2485145537
$ws.Range("AD2").Value = 6.3235353528
$ws.Range("AE2").Value = -23.8337505423
$ws.Range("AF2").Value = 689.5016373008
$ws.Range("AG2").Value = 11.7632194502
